$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.513003127792729
$ws.Range("D2").Value = 0.1445128998019207

# Row 3
$ws.Range("C3").Value = 0.4979572199634416
$ws.Range("D3").Value = 0.6234554242942036

# Row 4
$ws.Range("C4").Value = 1.817457031762177
$ws.Range("D4").Value = 0.08279391146765369
$ws.Range("G4").Value = "No"

# Row 5
$ws.Range("C5").Value = 0.01322798873740942
$ws.Range("D5").Value = 0.9895651270017631

# Row 6
$ws.Range("C6").Value = -1.096493018437365
$ws.Range("D6").Value = 0.2847210121154709

# Row 7
$ws.Range("C7").Value = 0.431002957956545
$ws.Range("D7").Value = 0.6706598422932974

# Row 8
$ws.Range("C8").Value = -1.516085941683251
$ws.Range("D8").Value = 0.1437355652488148

# Row 9
$ws.Range("C9").Value = 1.371819212740616
$ws.Range("D9").Value = 0.1839484597363366

# Row 10
$ws.Range("C10").Value = -0.6628848391257672
$ws.Range("D10").Value = 0.5142901991052837

# Row 11
$ws.Range("C11").Value = -1.716267355215658
$ws.Range("D11").Value = 0.1001629858674864

$wb.Save()
